$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct mis-spelled "Secteur" values in the Meteo table --------------
# "Torra di Murtella" -> "A Torra di Murtella"
$ws.Range("B18:B21").Value = "A Torra di Murtella"
# "Maffalcu" -> "Malfalcu"
$ws.Range("B30:B33").Value = "Malfalcu"

# --- Widen the "Secteur" column (B) so the longer names fit ---------------
$ws.Columns("B").ColumnWidth = 40.2857142857143

# --- Update the active selection to the corrected cells -------------------
$ws.Range("B18:B21").Select()
